{"js": "// Office.js (Word JavaScript API) script\n// Updates the date title and all two-digit multiplication problems in the\n// document body (1 title paragraph + 100 table-cell paragraphs, in document\n// order) to match the target revision.\n\nconst oldTexts = [\"2023-05-08 Monday\", \"82\u00d734=\", \"29\u00d770=\", \"78\u00d753=\", \"46\u00d732=\", \"41\u00d7100=\", \"39\u00d773=\", \"26\u00d772=\", \"86\u00d789=\", \"70\u00d771=\", \"62\u00d758=\", \"70\u00d736=\", \"90\u00d746=\", \"25\u00d793=\", \"19\u00d758=\", \"40\u00d795=\", \"39\u00d736=\", \"10\u00d780=\", \"98\u00d781=\", \"18\u00d779=\", \"65\u00d750=\", \"26\u00d751=\", \"11\u00d772=\", \"29\u00d768=\", \"92\u00d777=\", \"60\u00d777=\", \"12\u00d783=\", \"29\u00d733=\", \"51\u00d776=\", \"60\u00d712=\", \"38\u00d788=\", \"88\u00d759=\", \"40\u00d799=\", \"72\u00d777=\", \"32\u00d723=\", \"17\u00d770=\", \"78\u00d754=\", \"67\u00d771=\", \"55\u00d727=\", \"63\u00d735=\", \"62\u00d7100=\", \"97\u00d742=\", \"50\u00d730=\", \"83\u00d734=\", \"61\u00d774=\", \"100\u00d715=\", \"83\u00d780=\", \"23\u00d727=\", \"42\u00d745=\", \"63\u00d764=\", \"37\u00d774=\", \"49\u00d732=\", \"25\u00d780=\", \"88\u00d795=\", \"53\u00d749=\", \"57\u00d785=\", \"74\u00d760=\", \"94\u00d721=\", \"10\u00d795=\", \"52\u00d787=\", \"76\u00d796=\", \"32\u00d758=\", \"66\u00d759=\", \"84\u00d783=\", \"33\u00d787=\", \"43\u00d753=\", \"88\u00d791=\", \"83\u00d732=\", \"80\u00d787=\", \"52\u00d731=\", \"17\u00d764=\", \"76\u00d738=\", \"53\u00d714=\", \"47\u00d756=\", \"31\u00d774=\", \"15\u00d777=\", \"96\u00d755=\", \"49\u00d737=\", \"42\u00d794=\", \"88\u00d775=\", \"67\u00d741=\", \"35\u00d739=\", \"89\u00d764=\", \"43\u00d755=\", \"14\u00d748=\", \"13\u00d743=\", \"35\u00d733=\", \"82\u00d784=\", \"30\u00d723=\", \"85\u00d781=\", \"22\u00d772=\", \"46\u00d719=\", \"81\u00d749=\", \"83\u00d765=\", \"40\u00d750=\", \"40\u00d733=\", \"68\u00d721=\", \"72\u00d799=\", \"47\u00d739=\", \"97\u00d744=\", \"74\u00d783=\"];\nconst newTexts = [\"2023-05-09 Tuesday\", \"15\u00d796=\", \"82\u00d745=\", \"53\u00d724=\", \"63\u00d773=\", \"63\u00d781=\", \"76\u00d727=\", \"30\u00d797=\", \"92\u00d726=\", \"16\u00d7100=\", \"87\u00d776=\", \"32\u00d746=\", \"89\u00d724=\", \"99\u00d774=\", \"47\u00d715=\", \"78\u00d776=\", \"35\u00d781=\", \"28\u00d795=\", \"44\u00d757=\", \"41\u00d783=\", \"81\u00d741=\", \"34\u00d727=\", \"50\u00d736=\", \"41\u00d761=\", \"54\u00d728=\", \"65\u00d785=\", \"49\u00d733=\", \"24\u00d735=\", \"38\u00d759=\", \"33\u00d757=\", \"16\u00d757=\", \"24\u00d725=\", \"23\u00d781=\", \"12\u00d743=\", \"57\u00d751=\", \"88\u00d793=\", \"50\u00d789=\", \"17\u00d782=\", \"17\u00d717=\", \"61\u00d773=\", \"60\u00d783=\", \"73\u00d711=\", \"88\u00d770=\", \"89\u00d726=\", \"48\u00d717=\", \"93\u00d773=\", \"11\u00d797=\", \"10\u00d711=\", \"79\u00d746=\", \"24\u00d786=\", \"56\u00d749=\", \"58\u00d788=\", \"64\u00d785=\", \"88\u00d732=\", \"44\u00d759=\", \"86\u00d722=\", \"29\u00d739=\", \"94\u00d743=\", \"70\u00d794=\", \"28\u00d780=\", \"85\u00d728=\", \"56\u00d770=\", \"70\u00d784=\", \"27\u00d782=\", \"66\u00d793=\", \"24\u00d726=\", \"53\u00d769=\", \"59\u00d710=\", \"70\u00d745=\", \"83\u00d798=\", \"92\u00d734=\", \"40\u00d749=\", \"22\u00d723=\", \"56\u00d725=\", \"87\u00d712=\", \"89\u00d744=\", \"14\u00d745=\", \"14\u00d710=\", \"14\u00d774=\", \"88\u00d762=\", \"48\u00d714=\", \"92\u00d718=\", \"41\u00d770=\", \"84\u00d787=\", \"73\u00d752=\", \"33\u00d788=\", \"97\u00d729=\", \"10\u00d776=\", \"30\u00d766=\", \"47\u00d797=\", \"49\u00d774=\", \"36\u00d781=\", \"10\u00d763=\", \"95\u00d757=\", \"28\u00d732=\", \"98\u00d750=\", \"47\u00d734=\", \"38\u00d747=\", \"93\u00d739=\", \"35\u00d774=\", \"56\u00d721=\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== newTexts.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + newTexts.length + \" got \" + items.length\n  );\n}\n\nconst warnings = [];\nfor (let i = 0; i < items.length; i++) {\n  const current = items[i].text;\n  if (current === newTexts[i]) {\n    continue; // already up to date\n  }\n  if (current !== oldTexts[i]) {\n    // Paragraph text doesn't match what we expect to find before the edit;\n    // note it, but still apply the intended replacement by position.\n    warnings.push(\"paragraph \" + i + \": expected '\" + oldTexts[i] + \"' but found '\" + current + \"'\");\n  }\n  items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\nif (warnings.length > 0) {\n  return \"Completed with warnings: \" + warnings.join(\"; \");\n}\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Updates the date title and all two-digit multiplication problems to\n# match the target revision using Find/Replace (each old value is unique\n# in the document, so a simple ordered Find.Execute replace-all is safe).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-05-08 Monday\", \"2023-05-09 Tuesday\"),\n    @(\"82\u00d734=\", \"15\u00d796=\"),\n    @(\"29\u00d770=\", \"82\u00d745=\"),\n    @(\"78\u00d753=\", \"53\u00d724=\"),\n    @(\"46\u00d732=\", \"63\u00d773=\"),\n    @(\"41\u00d7100=\", \"63\u00d781=\"),\n    @(\"39\u00d773=\", \"76\u00d727=\"),\n    @(\"26\u00d772=\", \"30\u00d797=\"),\n    @(\"86\u00d789=\", \"92\u00d726=\"),\n    @(\"70\u00d771=\", \"16\u00d7100=\"),\n    @(\"62\u00d758=\", \"87\u00d776=\"),\n    @(\"70\u00d736=\", \"32\u00d746=\"),\n    @(\"90\u00d746=\", \"89\u00d724=\"),\n    @(\"25\u00d793=\", \"99\u00d774=\"),\n    @(\"19\u00d758=\", \"47\u00d715=\"),\n    @(\"40\u00d795=\", \"78\u00d776=\"),\n    @(\"39\u00d736=\", \"35\u00d781=\"),\n    @(\"10\u00d780=\", \"28\u00d795=\"),\n    @(\"98\u00d781=\", \"44\u00d757=\"),\n    @(\"18\u00d779=\", \"41\u00d783=\"),\n    @(\"65\u00d750=\", \"81\u00d741=\"),\n    @(\"26\u00d751=\", \"34\u00d727=\"),\n    @(\"11\u00d772=\", \"50\u00d736=\"),\n    @(\"29\u00d768=\", \"41\u00d761=\"),\n    @(\"92\u00d777=\", \"54\u00d728=\"),\n    @(\"60\u00d777=\", \"65\u00d785=\"),\n    @(\"12\u00d783=\", \"49\u00d733=\"),\n    @(\"29\u00d733=\", \"24\u00d735=\"),\n    @(\"51\u00d776=\", \"38\u00d759=\"),\n    @(\"60\u00d712=\", \"33\u00d757=\"),\n    @(\"38\u00d788=\", \"16\u00d757=\"),\n    @(\"88\u00d759=\", \"24\u00d725=\"),\n    @(\"40\u00d799=\", \"23\u00d781=\"),\n    @(\"72\u00d777=\", \"12\u00d743=\"),\n    @(\"32\u00d723=\", \"57\u00d751=\"),\n    @(\"17\u00d770=\", \"88\u00d793=\"),\n    @(\"78\u00d754=\", \"50\u00d789=\"),\n    @(\"67\u00d771=\", \"17\u00d782=\"),\n    @(\"55\u00d727=\", \"17\u00d717=\"),\n    @(\"63\u00d735=\", \"61\u00d773=\"),\n    @(\"62\u00d7100=\", \"60\u00d783=\"),\n    @(\"97\u00d742=\", \"73\u00d711=\"),\n    @(\"50\u00d730=\", \"88\u00d770=\"),\n    @(\"83\u00d734=\", \"89\u00d726=\"),\n    @(\"61\u00d774=\", \"48\u00d717=\"),\n    @(\"100\u00d715=\", \"93\u00d773=\"),\n    @(\"83\u00d780=\", \"11\u00d797=\"),\n    @(\"23\u00d727=\", \"10\u00d711=\"),\n    @(\"42\u00d745=\", \"79\u00d746=\"),\n    @(\"63\u00d764=\", \"24\u00d786=\"),\n    @(\"37\u00d774=\", \"56\u00d749=\"),\n    @(\"49\u00d732=\", \"58\u00d788=\"),\n    @(\"25\u00d780=\", \"64\u00d785=\"),\n    @(\"88\u00d795=\", \"88\u00d732=\"),\n    @(\"53\u00d749=\", \"44\u00d759=\"),\n    @(\"57\u00d785=\", \"86\u00d722=\"),\n    @(\"74\u00d760=\", \"29\u00d739=\"),\n    @(\"94\u00d721=\", \"94\u00d743=\"),\n    @(\"10\u00d795=\", \"70\u00d794=\"),\n    @(\"52\u00d787=\", \"28\u00d780=\"),\n    @(\"76\u00d796=\", \"85\u00d728=\"),\n    @(\"32\u00d758=\", \"56\u00d770=\"),\n    @(\"66\u00d759=\", \"70\u00d784=\"),\n    @(\"84\u00d783=\", \"27\u00d782=\"),\n    @(\"33\u00d787=\", \"66\u00d793=\"),\n    @(\"43\u00d753=\", \"24\u00d726=\"),\n    @(\"88\u00d791=\", \"53\u00d769=\"),\n    @(\"83\u00d732=\", \"59\u00d710=\"),\n    @(\"80\u00d787=\", \"70\u00d745=\"),\n    @(\"52\u00d731=\", \"83\u00d798=\"),\n    @(\"17\u00d764=\", \"92\u00d734=\"),\n    @(\"76\u00d738=\", \"40\u00d749=\"),\n    @(\"53\u00d714=\", \"22\u00d723=\"),\n    @(\"47\u00d756=\", \"56\u00d725=\"),\n    @(\"31\u00d774=\", \"87\u00d712=\"),\n    @(\"15\u00d777=\", \"89\u00d744=\"),\n    @(\"96\u00d755=\", \"14\u00d745=\"),\n    @(\"49\u00d737=\", \"14\u00d710=\"),\n    @(\"42\u00d794=\", \"14\u00d774=\"),\n    @(\"88\u00d775=\", \"88\u00d762=\"),\n    @(\"67\u00d741=\", \"48\u00d714=\"),\n    @(\"35\u00d739=\", \"92\u00d718=\"),\n    @(\"89\u00d764=\", \"41\u00d770=\"),\n    @(\"43\u00d755=\", \"84\u00d787=\"),\n    @(\"14\u00d748=\", \"73\u00d752=\"),\n    @(\"13\u00d743=\", \"33\u00d788=\"),\n    @(\"35\u00d733=\", \"97\u00d729=\"),\n    @(\"82\u00d784=\", \"10\u00d776=\"),\n    @(\"30\u00d723=\", \"30\u00d766=\"),\n    @(\"85\u00d781=\", \"47\u00d797=\"),\n    @(\"22\u00d772=\", \"49\u00d774=\"),\n    @(\"46\u00d719=\", \"36\u00d781=\"),\n    @(\"81\u00d749=\", \"10\u00d763=\"),\n    @(\"83\u00d765=\", \"95\u00d757=\"),\n    @(\"40\u00d750=\", \"28\u00d732=\"),\n    @(\"40\u00d733=\", \"98\u00d750=\"),\n    @(\"68\u00d721=\", \"47\u00d734=\"),\n    @(\"72\u00d799=\", \"38\u00d747=\"),\n    @(\"47\u00d739=\", \"93\u00d739=\"),\n    @(\"97\u00d744=\", \"35\u00d774=\"),\n    @(\"74\u00d783=\", \"56\u00d721=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: not found -> $oldText\"\n    }\n}\n\n"}
